$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A1").Value = 45436

$ws.Range("D10").Value = 3231
$ws.Range("D11").Value = 3913
$ws.Range("D12").Value = 3640
$ws.Range("D13").Value = 4459
$ws.Range("D14").Value = 5187
$ws.Range("D15").Value = 3822
$ws.Range("D16").Value = 6753
$ws.Range("D17").Value = 324

$ws.Range("D25").Value = 4050
$ws.Range("D26").Value = 7080
$ws.Range("D27").Value = 5460
$ws.Range("D28").Value = 7470
$ws.Range("D29").Value = 5100
$ws.Range("D30").Value = 8477
$ws.Range("D31").Value = 6734
$ws.Range("D32").Value = 8750
